# Corrección: Identificador de ruta.
# Updates "Tabla#1 - Error" (error-character log) rows 4-14 and
# "Tabla#2 - Tokens" (token log) rows 3-10 to reflect the corrected
# lexer output after fixing the "ruta" path identifier handling.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Tabla#1 - Error" ----
$wsErr = $wb.Worksheets.Item("Tabla#1 - Error")

# Row 4: now a short 2-row/no-column error whose character is a TAB
$wsErr.Range("B4").Value = 2
$wsErr.Range("C4").ClearContents()
$wsErr.Range("D4").Value = "`t"

# Row 5
$wsErr.Range("B5").Value = 6
$wsErr.Range("C5").Value = 22
$wsErr.Range("D5").Value = "@"

# Row 6
$wsErr.Range("B6").Value = 15
$wsErr.Range("C6").Value = 7
$wsErr.Range("D6").Value = "I"

# Row 7
$wsErr.Range("B7").Value = 22
$wsErr.Range("C7").Value = 7
$wsErr.Range("D7").Value = "I"

# Row 8
$wsErr.Range("B8").Value = 30
$wsErr.Range("C8").Value = 7
$wsErr.Range("D8").Value = "I"

# Row 9
$wsErr.Range("B9").Value = 37
$wsErr.Range("C9").Value = 2
$wsErr.Range("D9").Value = "@"

# Row 10
$wsErr.Range("B10").Value = 47
$wsErr.Range("C10").Value = 11
$wsErr.Range("D10").Value = "`$"

# Row 11
$wsErr.Range("B11").Value = 51
$wsErr.Range("C11").Value = 7
$wsErr.Range("D11").Value = "I"

# Row 12 (previously the last data row) now holds the old row-10 error
$wsErr.Range("A12").Value = 11
$wsErr.Range("B12").Value = 137
$wsErr.Range("C12").Value = 13
$wsErr.Range("D12").Value = [char]8221

# Row 13 (new): previously row-11 error
$wsErr.Range("A13").Value = 12
$wsErr.Range("B13").Value = 146
$wsErr.Range("C13").Value = 7
$wsErr.Range("D13").Value = "I"

# Row 14 (new): previously row-12 error
$wsErr.Range("A14").Value = 13
$wsErr.Range("B14").Value = 150
$wsErr.Range("C14").Value = 12
$wsErr.Range("D14").Value = "@"

# ---- Sheet 2: "Tabla#2 - Tokens" ----
$wsTok = $wb.Worksheets.Item("Tabla#2 - Tokens")

# Row 3
$wsTok.Range("B3").Value = "nombre"
$wsTok.Range("D3").Value = 3
$wsTok.Range("E3").Value = "nombre"

# Row 4
$wsTok.Range("B4").Value = "nombre"
$wsTok.Range("E4").Value = "nombre"

# Row 5
$wsTok.Range("B5").Value = "fin"
$wsTok.Range("E5").Value = "fin"

# Row 6
$wsTok.Range("B6").Value = "fin"
$wsTok.Range("E6").Value = "fin"

# Row 7
$wsTok.Range("B7").Value = "inicio"
$wsTok.Range("E7").Value = "inicio"

# Row 8
$wsTok.Range("B8").Value = "inicio"
$wsTok.Range("C8").Value = 10
$wsTok.Range("D8").Value = 7
$wsTok.Range("E8").Value = "inicio"

# Row 9
$wsTok.Range("B9").Value = "peso"
$wsTok.Range("C9").Value = 11
$wsTok.Range("E9").Value = "peso"

# Row 10
$wsTok.Range("B10").Value = "peso"
$wsTok.Range("D10").Value = 15
$wsTok.Range("E10").Value = "peso"
